$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert 4 new rows right above the "RC Command" section (old row 96) to host
# the four new Artisan commands: showCurve, showExtraCurve, showEvents,
# showBackgroundEvents. This pushes the existing "RC Command" block (and
# everything below it) down by 4 rows.
$ws.Rows("96:99").Insert()

# Match the row height used by the neighbouring keyboard-mode rows.
$ws.Rows("96:99").RowHeight = 13.8

# Row 96: showCurve(<name>,<bool>)
$ws.Cells.Item(96, 2).Value = "showCurve(<name>,<bool>)"
$ws.Cells.Item(96, 3).Value = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

# Row 97: showExtraCurve(<extra_device>,<curve>,<bool>)
$ws.Cells.Item(97, 2).Value = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$ws.Cells.Item(97, 3).Value = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

# Row 98: showEvents(<event_type>, <bool>)
$ws.Cells.Item(98, 2).Value = "showEvents(<event_type>, <bool>)"
$ws.Cells.Item(98, 3).Value = "shows/hides the events of <event_type> in [1,..,5]"

# Row 99: showBackgroundEvents(<bool>)
$ws.Cells.Item(99, 2).Value = "showBackgroundEvents(<bool>)"
$ws.Cells.Item(99, 3).Value = "shows/hides the events of the background profile"

# Update the selection to mirror the new active cell after the edit, and
# make sure "Commands" remains the active sheet/tab.
$ws1 = $wb.Worksheets.Item("Sliders")
$ws1.Range("B6").Select()

$ws.Select()
$ws.Range("C97").Select()
